# Regenerate orders with updated distance/size codes.
# The experiment's distance and size condition labels changed:
#   D64 -> D69, D80 -> D86, D51 -> D55 (distance codes)
#   S30 -> S31 (size code)
# These codes appear embedded inside many shared strings (condition
# names like "Face16_D64_S20", image filenames like
# "Face16_D64_S20_l.png", and the standalone lookup values "D64"/"S30"),
# so apply the substitutions as workbook-wide text replacements across
# every used cell, mirroring Excel's Find & Replace ("Replace All").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$used.Replace("D64", "D69")
$used.Replace("D80", "D86")
$used.Replace("D51", "D55")
$used.Replace("S30", "S31")
